# Actualización desde MV -datos-
# Updates quarterly export figures (CUCI) for rows 74 and 75.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 74
$ws.Range("B74").Value = 22232
$ws.Range("C74").Value = 5205
$ws.Range("D74").Value = 1488
$ws.Range("E74").Value = 3119
$ws.Range("I74").Value = 57
$ws.Range("N74").Value = 8689
$ws.Range("O74").Value = 7695
$ws.Range("Q74").Value = 302
$ws.Range("R74").Value = 106
$ws.Range("T74").Value = 156
$ws.Range("U74").Value = 123
$ws.Range("X74").Value = 820
$ws.Range("Z74").Value = 444
$ws.Range("AC74").Value = 6124
$ws.Range("AD74").Value = 5439
$ws.Range("AF74").Value = 128
$ws.Range("AI74").Value = 201

# Row 75
$ws.Range("B75").Value = 23234
$ws.Range("C75").Value = 3793
$ws.Range("D75").Value = 1445
$ws.Range("E75").Value = 1563
$ws.Range("G75").Value = 386
$ws.Range("H75").Value = 83
$ws.Range("I75").Value = 108
$ws.Range("N75").Value = 10194
$ws.Range("O75").Value = 8915
$ws.Range("Q75").Value = 385
$ws.Range("T75").Value = 136
$ws.Range("U75").Value = 112
$ws.Range("X75").Value = 1033
$ws.Range("Z75").Value = 633
$ws.Range("AC75").Value = 6702
$ws.Range("AD75").Value = 5971
$ws.Range("AG75").Value = 94
$ws.Range("AH75").Value = 44
$ws.Range("AJ75").Value = 342
$ws.Range("AK75").Value = 135
$ws.Range("AL75").Value = 207
$ws.Range("AM75").Value = 149
$ws.Range("AN75").Value = 43
$ws.Range("AO75").Value = 105
